# Auto-generated: append new sensor log rows to PIR, Humidity, Temperature, and Proximity sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(242, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(242, 2).Value = "'" + '15:12:16'
$ws.Cells.Item(242, 3).Value = "'" + '15:00'
$ws.Cells.Item(242, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(242, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(242, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(243, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(243, 2).Value = "'" + '15:12:18'
$ws.Cells.Item(243, 3).Value = "'" + '15:00'
$ws.Cells.Item(243, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(243, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(243, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(244, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(244, 2).Value = "'" + '15:12:22'
$ws.Cells.Item(244, 3).Value = "'" + '15:00'
$ws.Cells.Item(244, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(244, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(244, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(245, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(245, 2).Value = "'" + '15:12:28'
$ws.Cells.Item(245, 3).Value = "'" + '15:00'
$ws.Cells.Item(245, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(245, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(245, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(246, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(246, 2).Value = "'" + '15:12:32'
$ws.Cells.Item(246, 3).Value = "'" + '15:00'
$ws.Cells.Item(246, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(246, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(246, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(247, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(247, 2).Value = "'" + '15:12:37'
$ws.Cells.Item(247, 3).Value = "'" + '15:00'
$ws.Cells.Item(247, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(247, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(247, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(248, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(248, 2).Value = "'" + '15:12:42'
$ws.Cells.Item(248, 3).Value = "'" + '15:00'
$ws.Cells.Item(248, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(248, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(248, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(249, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(249, 2).Value = "'" + '15:12:48'
$ws.Cells.Item(249, 3).Value = "'" + '15:00'
$ws.Cells.Item(249, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(249, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(249, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(250, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(250, 2).Value = "'" + '15:12:52'
$ws.Cells.Item(250, 3).Value = "'" + '15:00'
$ws.Cells.Item(250, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(250, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(250, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(251, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(251, 2).Value = "'" + '15:12:57'
$ws.Cells.Item(251, 3).Value = "'" + '15:00'
$ws.Cells.Item(251, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(251, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(251, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(252, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(252, 2).Value = "'" + '15:13:02'
$ws.Cells.Item(252, 3).Value = "'" + '15:00'
$ws.Cells.Item(252, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(252, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(252, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(253, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(253, 2).Value = "'" + '15:13:08'
$ws.Cells.Item(253, 3).Value = "'" + '15:00'
$ws.Cells.Item(253, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(253, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(253, 6).Value = "'" + 'Inactive'
$ws.Cells.Item(254, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(254, 2).Value = "'" + '15:13:12'
$ws.Cells.Item(254, 3).Value = "'" + '15:00'
$ws.Cells.Item(254, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(254, 5).Value = "'" + 'No Motion'
$ws.Cells.Item(254, 6).Value = "'" + 'Inactive'

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(231, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(231, 2).Value = "'" + '15:12:17'
$ws.Cells.Item(231, 3).Value = "'" + '15:00'
$ws.Cells.Item(231, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(231, 5).Value = "'" + '88.3%'
$ws.Cells.Item(231, 6).Value = "'" + 'Active'
$ws.Cells.Item(232, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(232, 2).Value = "'" + '15:12:19'
$ws.Cells.Item(232, 3).Value = "'" + '15:00'
$ws.Cells.Item(232, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(232, 5).Value = "'" + '88.3%'
$ws.Cells.Item(232, 6).Value = "'" + 'Active'
$ws.Cells.Item(233, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(233, 2).Value = "'" + '15:12:23'
$ws.Cells.Item(233, 3).Value = "'" + '15:00'
$ws.Cells.Item(233, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(233, 5).Value = "'" + '87.4%'
$ws.Cells.Item(233, 6).Value = "'" + 'Active'
$ws.Cells.Item(234, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(234, 2).Value = "'" + '15:12:26'
$ws.Cells.Item(234, 3).Value = "'" + '15:00'
$ws.Cells.Item(234, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(234, 5).Value = "'" + '88.3%'
$ws.Cells.Item(234, 6).Value = "'" + 'Active'
$ws.Cells.Item(235, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(235, 2).Value = "'" + '15:12:30'
$ws.Cells.Item(235, 3).Value = "'" + '15:00'
$ws.Cells.Item(235, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(235, 5).Value = "'" + '88.4%'
$ws.Cells.Item(235, 6).Value = "'" + 'Active'
$ws.Cells.Item(236, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(236, 2).Value = "'" + '15:12:34'
$ws.Cells.Item(236, 3).Value = "'" + '15:00'
$ws.Cells.Item(236, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(236, 5).Value = "'" + '87.4%'
$ws.Cells.Item(236, 6).Value = "'" + 'Active'
$ws.Cells.Item(237, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(237, 2).Value = "'" + '15:12:38'
$ws.Cells.Item(237, 3).Value = "'" + '15:00'
$ws.Cells.Item(237, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(237, 5).Value = "'" + '88.3%'
$ws.Cells.Item(237, 6).Value = "'" + 'Active'
$ws.Cells.Item(238, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(238, 2).Value = "'" + '15:12:43'
$ws.Cells.Item(238, 3).Value = "'" + '15:00'
$ws.Cells.Item(238, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(238, 5).Value = "'" + '87.5%'
$ws.Cells.Item(238, 6).Value = "'" + 'Active'
$ws.Cells.Item(239, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(239, 2).Value = "'" + '15:12:47'
$ws.Cells.Item(239, 3).Value = "'" + '15:00'
$ws.Cells.Item(239, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(239, 5).Value = "'" + '88.3%'
$ws.Cells.Item(239, 6).Value = "'" + 'Active'
$ws.Cells.Item(240, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(240, 2).Value = "'" + '15:12:51'
$ws.Cells.Item(240, 3).Value = "'" + '15:00'
$ws.Cells.Item(240, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(240, 5).Value = "'" + '88.4%'
$ws.Cells.Item(240, 6).Value = "'" + 'Active'
$ws.Cells.Item(241, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(241, 2).Value = "'" + '15:12:55'
$ws.Cells.Item(241, 3).Value = "'" + '15:00'
$ws.Cells.Item(241, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(241, 5).Value = "'" + '87.5%'
$ws.Cells.Item(241, 6).Value = "'" + 'Active'
$ws.Cells.Item(242, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(242, 2).Value = "'" + '15:13:07'
$ws.Cells.Item(242, 3).Value = "'" + '15:00'
$ws.Cells.Item(242, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(242, 5).Value = "'" + '88.3%'
$ws.Cells.Item(242, 6).Value = "'" + 'Active'
$ws.Cells.Item(243, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(243, 2).Value = "'" + '15:13:11'
$ws.Cells.Item(243, 3).Value = "'" + '15:00'
$ws.Cells.Item(243, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(243, 5).Value = "'" + '88.4%'
$ws.Cells.Item(243, 6).Value = "'" + 'Active'
$ws.Cells.Item(244, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(244, 2).Value = "'" + '15:13:15'
$ws.Cells.Item(244, 3).Value = "'" + '15:00'
$ws.Cells.Item(244, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(244, 5).Value = "'" + '87.5%'
$ws.Cells.Item(244, 6).Value = "'" + 'Active'

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(231, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(231, 2).Value = "'" + '15:12:17'
$ws.Cells.Item(231, 3).Value = "'" + '15:00'
$ws.Cells.Item(231, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(231, 5).Value = "'" + '23.0C'
$ws.Cells.Item(231, 6).Value = "'" + 'Active'
$ws.Cells.Item(232, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(232, 2).Value = "'" + '15:12:19'
$ws.Cells.Item(232, 3).Value = "'" + '15:00'
$ws.Cells.Item(232, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(232, 5).Value = "'" + '22.9C'
$ws.Cells.Item(232, 6).Value = "'" + 'Active'
$ws.Cells.Item(233, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(233, 2).Value = "'" + '15:12:23'
$ws.Cells.Item(233, 3).Value = "'" + '15:00'
$ws.Cells.Item(233, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(233, 5).Value = "'" + '22.9C'
$ws.Cells.Item(233, 6).Value = "'" + 'Active'
$ws.Cells.Item(234, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(234, 2).Value = "'" + '15:12:27'
$ws.Cells.Item(234, 3).Value = "'" + '15:00'
$ws.Cells.Item(234, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(234, 5).Value = "'" + '22.9C'
$ws.Cells.Item(234, 6).Value = "'" + 'Active'
$ws.Cells.Item(235, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(235, 2).Value = "'" + '15:12:31'
$ws.Cells.Item(235, 3).Value = "'" + '15:00'
$ws.Cells.Item(235, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(235, 5).Value = "'" + '23.0C'
$ws.Cells.Item(235, 6).Value = "'" + 'Active'
$ws.Cells.Item(236, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(236, 2).Value = "'" + '15:12:35'
$ws.Cells.Item(236, 3).Value = "'" + '15:00'
$ws.Cells.Item(236, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(236, 5).Value = "'" + '23.0C'
$ws.Cells.Item(236, 6).Value = "'" + 'Active'
$ws.Cells.Item(237, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(237, 2).Value = "'" + '15:12:39'
$ws.Cells.Item(237, 3).Value = "'" + '15:00'
$ws.Cells.Item(237, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(237, 5).Value = "'" + '22.9C'
$ws.Cells.Item(237, 6).Value = "'" + 'Active'
$ws.Cells.Item(238, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(238, 2).Value = "'" + '15:12:43'
$ws.Cells.Item(238, 3).Value = "'" + '15:00'
$ws.Cells.Item(238, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(238, 5).Value = "'" + '23.0C'
$ws.Cells.Item(238, 6).Value = "'" + 'Active'
$ws.Cells.Item(239, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(239, 2).Value = "'" + '15:12:47'
$ws.Cells.Item(239, 3).Value = "'" + '15:00'
$ws.Cells.Item(239, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(239, 5).Value = "'" + '22.9C'
$ws.Cells.Item(239, 6).Value = "'" + 'Active'
$ws.Cells.Item(240, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(240, 2).Value = "'" + '15:12:51'
$ws.Cells.Item(240, 3).Value = "'" + '15:00'
$ws.Cells.Item(240, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(240, 5).Value = "'" + '23.0C'
$ws.Cells.Item(240, 6).Value = "'" + 'Active'
$ws.Cells.Item(241, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(241, 2).Value = "'" + '15:12:55'
$ws.Cells.Item(241, 3).Value = "'" + '15:00'
$ws.Cells.Item(241, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(241, 5).Value = "'" + '23.0C'
$ws.Cells.Item(241, 6).Value = "'" + 'Active'
$ws.Cells.Item(242, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(242, 2).Value = "'" + '15:13:07'
$ws.Cells.Item(242, 3).Value = "'" + '15:00'
$ws.Cells.Item(242, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(242, 5).Value = "'" + '22.9C'
$ws.Cells.Item(242, 6).Value = "'" + 'Active'
$ws.Cells.Item(243, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(243, 2).Value = "'" + '15:13:11'
$ws.Cells.Item(243, 3).Value = "'" + '15:00'
$ws.Cells.Item(243, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(243, 5).Value = "'" + '23.0C'
$ws.Cells.Item(243, 6).Value = "'" + 'Active'
$ws.Cells.Item(244, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(244, 2).Value = "'" + '15:13:16'
$ws.Cells.Item(244, 3).Value = "'" + '15:00'
$ws.Cells.Item(244, 4).Value = "'" + 'Bathroom'
$ws.Cells.Item(244, 5).Value = "'" + '22.9C'
$ws.Cells.Item(244, 6).Value = "'" + 'Active'

$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(32, 1).Value = "'" + '2026-01-28'
$ws.Cells.Item(32, 2).Value = "'" + '15:13:14'
$ws.Cells.Item(32, 3).Value = "'" + '15:00'
$ws.Cells.Item(32, 4).Value = "'" + 'Bedroom Door'
$ws.Cells.Item(32, 5).Value = "'" + 'Clear'
$ws.Cells.Item(32, 6).Value = "'" + 'Inactive'
